$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '60.027.56'
$ws.Range('D3').Value = '3.186.88'
$ws.Range('E3').Value = '  +1.09%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'536.40"
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('D6').Value = "'145.07"
$ws.Range('E6').Value = '  +4.03%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = "'0.533"
$ws.Range('E8').Value = '  -1.12%  '
$ws.Range('D9').Value = "'7.33"
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('E11').Value = '  -1.76%  '
$ws.Range('D12').Value = '3.737.65'
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('E13').Value = '  -2.84%  '
$ws.Range('D14').Value = "'25.80"
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').Value = "'0.0000172"
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').Value = '60.035.01'
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('D17').Value = '3.190.77'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('E18').Value = '  -0.38%  '
$ws.Range('D19').Value = "'13.23"
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('D20').Value = "'8.18"
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').Value = "'368.80"
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = "'0.522"
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('D24').Value = "'69.46"
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  +1.22%  '
$ws.Range('D26').Value = "'8.60"
$ws.Range('E26').Value = '  +3.41%  '
$ws.Range('D27').Value = "'0.983"
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('D28').Value = '0.0₃0874'
$ws.Range('E28').Value = '  +1.09%  '
$ws.Range('D29').Value = "'22.47"
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('E32').Value = '  +2.88%  '
$ws.Range('E33').Value = '  +5.00%  '
$ws.Range('E34').Value = '  +3.02%  '
$ws.Range('D35').Value = "'157.72"
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('E36').Value = '  +1.83%  '
$ws.Range('D37').Value = "'26.42"
$ws.Range('E37').Value = '  +5.88%  '
$ws.Range('D38').Value = '2.789.93'
$ws.Range('E38').Value = '  +5.77%  '
$ws.Range('E39').Value = '  +3.48%  '
$ws.Range('D40').Value = "'0.0308"
$ws.Range('E40').Value = '  +7.64%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').Value = "'39.81"
$ws.Range('E43').Value = '  +2.00%  '
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('E45').Value = '  +1.67%  '
$ws.Range('D46').Value = '3.228.27'
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('D47').Value = "'0.983"
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').Value = "'0.796"
$ws.Range('E49').Value = '  +5.96%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = "'20.58"
$ws.Range('E50').Value = '  +3.06%  '
